# Update the "HPEbP" sheet's natural gas reforming efficiency formula
# (cell B3) to drop the "+46" term, per the commit:
# "Update ELF and HPEbP with customizations from analysis repo to
# accomodate BAU hydrogen production"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HPEbP")

$ws.Range("B3").Formula = "=118/(162+2)"

$wb.Application.Calculate()
